# Add season-record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new labels in AD, AE, AF. Clone the formatting of
# the existing header cell (bold / centered / bordered) onto each new
# header cell before writing its text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-47: every player on the roster shares the team's season
# record: 69 wins, 93 losses, 0 ties.
for ($row = 2; $row -le 47; $row++) {
    $ws.Cells.Item($row, 30).Value = 69
    $ws.Cells.Item($row, 31).Value = 93
    $ws.Cells.Item($row, 32).Value = 0
}
